# Adding hashtable to parameterize
#
# - test_suite: the "openAccountTest" run-mode flips from N -> Y.
# - AddCustomerTest: a new "runmode" column (E) is added, set to "Y"
#   for every data row, so the sheet can be driven from the same
#   hashtable-based parameterization as the other suites.
# - Selection/active-cell bookmarks move around, and AddCustomerTest
#   becomes the active tab.

$wb = $excel.ActiveWorkbook

$wsSuite  = $wb.Worksheets.Item("test_suite")
$wsAdd    = $wb.Worksheets.Item("AddCustomerTest")
$wsOpen   = $wb.Worksheets.Item("OpenAccountTest")

# --- test_suite: flip openAccountTest's runmode from N to Y ---
$wsSuite.Range("B4").Value = "Y"

# --- AddCustomerTest: add the runmode column (E) ---
$wsAdd.Range("E1").Value = "runmode"
$wsAdd.Range("E2").Value = "Y"
$wsAdd.Range("E3").Value = "Y"
$wsAdd.Range("E4").Value = "Y"
$wsAdd.Range("E5").Value = "Y"

# --- Update each sheet's stored selection ---
$wsSuite.Range("B2").Select() | Out-Null
$wsAdd.Range("E5").Select() | Out-Null
$wsOpen.Range("C1").Select() | Out-Null

# --- AddCustomerTest becomes the active sheet/tab ---
$wsAdd.Activate() | Out-Null
